# Apply "refine figures and ch1" font fix: the four small caption text
# boxes ("白盒模型", "预测行为抽取", "可解释模型", "可解释决策路径覆盖度")
# get an explicit SimSun (宋体) Latin/East-Asian typeface on their run
# text, matching what PowerPoint writes when a Chinese font is applied
# from the Font dropdown.

$p = $ppt.ActivePresentation

$targetTexts = @("白盒模型", "预测行为抽取", "可解释模型", "可解释决策路径覆盖度")
$fontName = "宋体"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)

        if (-not $shape.HasTextFrame) {
            continue
        }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) {
            continue
        }

        $tr = $tf.TextRange
        if ($targetTexts -contains $tr.Text) {
            $tr.Font.Name = $fontName
            $tr.Font.NameFarEast = $fontName
        }
    }
}
